$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 2023 first data delivery - add the 2021 row of SAV acreage data to the
# existing table (row 39 was previously a blank placeholder row).
$ws.Range("A39").Value = 2021
$ws.Range("B39").Value = 19173
$ws.Range("C39").Value = "NA"
$ws.Range("D39").Value = "NA"
$ws.Range("E39").Value = 16132
$ws.Range("F39").Value = "NA"

# Reflect the scrolled/selected view used while reviewing the new data.
$ws.Activate()
$ws.Range("F42:F43").Select() | Out-Null
